$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11:F11").Copy($ws.Range("A12:F12"))
$ws.Range("A12").Value = 159
$ws.Range("B12").Value = 0.00000126721215565
$ws.Range("C12").Value = 15868943608
$ws.Range("F12").Value = "c86896c2d9ac90c9382981b8744282aafd76d861 "
$ws.Range("D12").Value = "Divide et impera"
$ws.Range("E12").Value = 42864

$co1 = $ws.ChartObjects(1)
$s1 = $co1.Chart.SeriesCollection().Item(1)
$arr = @(16321447886,16254382549,16313521178,15952716079,15949435607,15951842731,15968485694,16007210792,15962151054,15967660563,15868943608)
$s1.Values = $arr
$s1.Formula = "=SERIES(Tabelle1!`$C`$1,,Tabelle1!`$C`$2:`$C`$12,1)"
